$d = $word.ActiveDocument

# --- First paragraph: **ID__AFFARS_pgi_5305_topic_4__ID** / " " ---
$p1 = $d.Paragraphs.Item(1)

# Add a paragraph border (top/left/bottom/right) with 5pt space-from-text,
# no explicit line (matches <w:pBdr><w:top w:space="5"/>...).
$borders = $p1.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromRight = 5

# Change left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Format.LeftIndent = 11.25

# Update the placeholder id text (keeps this run's own rPr).
$start = $p1.Range.Start
$run1Range = $d.Range($start, $start + 35)
$run1Range.Find.Execute("pgi_5305_topic_4", $true, $false, $false, $false, $false, $true, 1, $false, "AFMC_PGI_5305", 2)

# Remove the trailing " " run that followed the id text.
$p1b = $d.Paragraphs.Item(1)
$pStart = $p1b.Range.Start
$fullText = $p1b.Range.Text
$spaceIdx = $fullText.Length - 2
$spaceRange = $d.Range($pStart + $spaceIdx, $pStart + $spaceIdx + 1)
$spaceRange.Delete()
